$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1675
$ws.Range("J17").Value = 1675
$ws.Range("L17").Value = 5025
$ws.Range("N17").Value = -5361

$ws.Range("H69").Value = 4571.2
$ws.Range("I69").Value = 4542
$ws.Range("K69").Value = 13626
$ws.Range("M69").Value = -12752

$ws.Range("H72").Value = 4571.2
$ws.Range("I72").Value = 4542
$ws.Range("K72").Value = 40878
$ws.Range("M72").Value = -36510

$ws.Range("H98").Value = 2962.8845
$ws.Range("I98").Value = 2820.182
$ws.Range("K98").Value = 2820.182
$ws.Range("M98").Value = -1322.182

$ws.Range("H113").Value = 3695.7144
$ws.Range("I113").Value = 4087.5
$ws.Range("J113").Value = 3173.3333
$ws.Range("K113").Value = 4087.5
$ws.Range("L113").Value = 3173.3333
$ws.Range("M113").Value = -833.5
$ws.Range("N113").Value = -9681.3333

$ws.Range("H122").Value = 2962.8845
$ws.Range("I122").Value = 2820.182
$ws.Range("K122").Value = 8460.545999999998
$ws.Range("M122").Value = -6010.545999999998

$ws.Range("H133").Value = 28200
$ws.Range("J133").Value = 28200
$ws.Range("L133").Value = 28200
$ws.Range("N133").Value = -38320

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6851.89
$ws.Range("I32").Value = 5949.9453
$ws.Range("J32").Value = 14167.667
$ws.Range("K32").Value = 5949.9453
$ws.Range("L32").Value = 14167.667
$ws.Range("M32").Value = -5662.9453
$ws.Range("N32").Value = -14741.667

$ws.Range("H97").Value = 777.7406999999999
$ws.Range("I97").Value = 715.86365
$ws.Range("J97").Value = 1050
$ws.Range("K97").Value = 715.86365
$ws.Range("L97").Value = 1050
$ws.Range("M97").Value = -219.86365
$ws.Range("N97").Value = -2042

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 18003.688
$ws.Range("I86").Value = 1948.2858
$ws.Range("J86").Value = 48654.91
$ws.Range("K86").Value = 1948.2858
$ws.Range("L86").Value = 48654.91
$ws.Range("M86").Value = -825.2858000000001
$ws.Range("N86").Value = -50900.91

$ws.Range("H89").Value = 18003.688
$ws.Range("I89").Value = 1948.2858
$ws.Range("J89").Value = 48654.91
$ws.Range("K89").Value = 9741.429
$ws.Range("L89").Value = 243274.55
$ws.Range("M89").Value = -4125.429
$ws.Range("N89").Value = -254506.55

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2306.535
$ws.Range("I31").Value = 1716.2703
$ws.Range("J31").Value = 5946.5
$ws.Range("K31").Value = 1716.2703
$ws.Range("L31").Value = 5946.5
$ws.Range("M31").Value = -1421.2703
$ws.Range("N31").Value = -6536.5

$ws.Range("H34").Value = 2306.535
$ws.Range("I34").Value = 1716.2703
$ws.Range("J34").Value = 5946.5
$ws.Range("K34").Value = 1716.2703
$ws.Range("L34").Value = 5946.5
$ws.Range("M34").Value = -1514.2703
$ws.Range("N34").Value = -6350.5

$ws.Range("H134").Value = 2113.7144
$ws.Range("I134").Value = 907.53845
$ws.Range("K134").Value = 2722.61535
$ws.Range("M134").Value = -187.61535

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 114.90909
$ws.Range("I23").Value = 79
$ws.Range("J23").Value = 122.888885
$ws.Range("K23").Value = 237
$ws.Range("L23").Value = 368.666655
$ws.Range("M23").Value = -2
$ws.Range("N23").Value = -838.666655

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3434.4358
$ws.Range("I132").Value = 3237.8
$ws.Range("J132").Value = 3785.5715
$ws.Range("K132").Value = 9713.400000000001
$ws.Range("L132").Value = 11356.7145
$ws.Range("M132").Value = -7183.400000000001
$ws.Range("N132").Value = -16416.7145

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4168242.5
$ws.Range("I7").Value = 6667621
$ws.Range("K7").Value = 6667621
$ws.Range("M7").Value = -6667509

$ws.Range("H68").Value = 1950
$ws.Range("I68").Value = 1056.25
$ws.Range("J68").Value = 5525
$ws.Range("K68").Value = 1056.25
$ws.Range("L68").Value = 5525
$ws.Range("M68").Value = -307.25
$ws.Range("N68").Value = -7023

$ws.Range("H71").Value = 1950
$ws.Range("I71").Value = 1056.25
$ws.Range("J71").Value = 5525
$ws.Range("K71").Value = 5281.25
$ws.Range("L71").Value = 27625
$ws.Range("M71").Value = -1537.25
$ws.Range("N71").Value = -35113

$ws.Range("H82").Value = 2704.1333
$ws.Range("I82").Value = 2039
$ws.Range("J82").Value = 3464.2856
$ws.Range("K82").Value = 2039
$ws.Range("L82").Value = 3464.2856
$ws.Range("M82").Value = -1678
$ws.Range("N82").Value = -4186.2856

$ws.Range("H85").Value = 2704.1333
$ws.Range("I85").Value = 2039
$ws.Range("J85").Value = 3464.2856
$ws.Range("K85").Value = 2039
$ws.Range("L85").Value = 3464.2856
$ws.Range("M85").Value = -791
$ws.Range("N85").Value = -5960.2856

$ws.Range("H126").Value = 4168242.5
$ws.Range("I126").Value = 6667621
$ws.Range("K126").Value = 20002863
$ws.Range("M126").Value = -20000393

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4818.1
$ws.Range("I62").Value = 4666.6665
$ws.Range("J62").Value = 4883
$ws.Range("K62").Value = 4666.6665
$ws.Range("L62").Value = 4883
$ws.Range("M62").Value = -4042.6665
$ws.Range("N62").Value = -6131

$ws.Range("H65").Value = 4818.1
$ws.Range("I65").Value = 4666.6665
$ws.Range("J65").Value = 4883
$ws.Range("K65").Value = 23333.3325
$ws.Range("L65").Value = 24415
$ws.Range("M65").Value = -20213.3325
$ws.Range("N65").Value = -30655

$ws.Range("H94").Value = 29200
$ws.Range("J94").Value = 29200
$ws.Range("L94").Value = 29200
$ws.Range("N94").Value = -31002

$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()

$ws.Range("H103").Value = 25638.477
$ws.Range("J103").Value = 25638.477
$ws.Range("L103").Value = 25638.477
$ws.Range("N103").Value = -27982.477

$ws.Range("H106").Value = 29297.125
$ws.Range("J106").Value = 29297.125
$ws.Range("L106").Value = 29297.125
$ws.Range("N106").Value = -31821.125

$ws.Range("H126").Value = 3227728.5
$ws.Range("I126").Value = 1185.3182
$ws.Range("J126").Value = 11114834
$ws.Range("K126").Value = 3555.9546
$ws.Range("L126").Value = 33344502
$ws.Range("M126").Value = -1085.9546
$ws.Range("N126").Value = -33349442
